$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(-0.0227839658534362, 0.6421757934121076, 0.7299546130156783, 0.8543738133953301, 0.8732649444381481, 23)
    3  = @(0.7639704308143828, 1.10269471751686, 2.612846386011487, 1.616430136446202, 1.458020551013975, 22)
    4  = @(0.6235520873978828, 1.3761391622945, 4.053773020290964, 2.013398375953195, 1.96168386987362, 21)
    5  = @(0.7192758613889139, 0.8905910875909095, 1.471357591919539, 1.212995297566952, 1.002102960958012, 20)
    6  = @(0.541769264226673, 0.8089314641762144, 1.078722909040149, 1.038615862116572, 0.9104019409704378, 19)
    7  = @(0.3378155251717742, 0.6939632015921933, 0.7224693225991232, 0.8499819542785148, 0.8025803930388514, 18)
    8  = @(0.276410312919357, 0.7156251782359094, 0.7145716310038928, 0.8453233884164644, 0.8234406660685688, 17)
    9  = @(0.553324246477921, 0.722017931190802, 0.7182827940224992, 0.8475156600455824, 0.6705073294552746, 12)
    10 = @(0.4056188661620858, 0.7064691664491064, 0.7927601511797749, 0.8903707942086684, 0.856118606089144, 7)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 2]
    }
}
